# Update app architecture sketches
#
# - Remove the "yelb-cache" can shape and its connector to yelb-appserver
# - Remove the "redis" label and its "5" step-number oval
# - Relocate the remaining "4" step-number oval down near the Postgres/db area

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# EMU -> Points helper. PowerPoint's Shape.Left/Top are single-precision
# points (1 pt = 12700 EMU); nudge by half an EMU so the float32 round-trip
# truncates back to the exact EMU value instead of landing one EMU short.
function ConvertTo-Points($emu) {
    return ($emu / 12700.0) + (0.5 / 12700.0)
}

function Set-ShapePositionEmu($shape, $xEmu, $yEmu) {
    $shape.Left = ConvertTo-Points $xEmu
    $shape.Top = ConvertTo-Points $yEmu
}

# Remove the "yelb-cache" can and the elbow connector that linked it to
# the appserver rectangle.
$s.Shapes.Item("Elbow Connector 9").Delete()
$s.Shapes.Item("Can 5").Delete()

# Reposition the "4" step oval that used to sit beside yelb-cache, moving
# it down near the Postgres / yelb-db area.
Set-ShapePositionEmu $s.Shapes.Item("Oval 52") 6979785 6177064

# Remove the "redis" callout rectangle and the "5" step oval that went
# with the deleted cache shape.
$s.Shapes.Item("Rounded Rectangle 38").Delete()
$s.Shapes.Item("Oval 41").Delete()
